$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelling "as_ravelling" -> "as_raveling" (row 4, defects_type column J)
$ws.Range("J4").Value = "as_raveling"

# The three defect-measurement columns (defects_width, defects_length, defects_vol)
# are collapsed into a single "defects_dimension" column. Drop the first two
# (K:L) so the former "defects_vol" column slides left into K, then overwrite
# its header/values with the new single dimension column.
$ws.Range("K:L").Delete()

$ws.Range("K1").Value = "defects_dimension"
$ws.Range("K2").Value = 20
$ws.Range("K3").Value = 20
$ws.Range("K4").Value = 20
$ws.Range("K5").Value = 20

# Update the sheet view to match the new layout (the active/selected cell
# anchor shifted now that two columns were removed).
$ws.Range("J5").Select()
